$d = $word.ActiveDocument

# Paragraph 1 currently holds the old title/subtitle text ("Reflections
# During Advent," Part One / "Searching for Christ ====...") spread over a
# <w:br/> and several runs. Replace it wholesale with a single "Title"
# styled paragraph containing "Dorothy Day" split into the same 3 runs the
# target markup uses ("Dorothy", " ", "Day").
$p1 = $d.Paragraphs.Item(1)
$titleXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
            '<w:pPr><w:pStyle w:val="Title"/></w:pPr>' +
            '<w:r><w:t xml:space="preserve">Dorothy</w:t></w:r>' +
            '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
            '<w:r><w:t xml:space="preserve">Day</w:t></w:r>' +
            '</w:p>'
$p1.Range.InsertXML($titleXml)

# Paragraph 2 ("By Dorothy Day") is no longer needed now that the byline is
# folded into the title paragraph above, so drop it (and its paragraph
# mark) entirely.
$p2 = $d.Paragraphs.Item(2)
$p2.Range.Delete()
